$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2, D2, E2 cleared; C2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.7011946562634641
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 values updated
$ws.Range("B3").Value = 6.0617727471557927
$ws.Range("C3").Value = 6.3386363225090436
$ws.Range("D3").Value = 7.798938069828
$ws.Range("E3").Value = 2.8392023425895161

# Update selection to match new range
$ws.Range("B1:E3").Select()
